# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff & handback datetime
# stamps for the 7e08dc8d-... entry across the Overview, zh-cn and de-de
# sheets to reflect the freshly generated report timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the 7e08dc8d-... row
$wsOverview.Range("G4").Value = "2016-08-24 22:46:12"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 7e08dc8d-... row
$wsZhCn.Range("H4").Value = "2016-08-24 22:46:08"
$wsZhCn.Range("K4").Value = "2016-08-24 22:46:25"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 7e08dc8d-... row
$wsDeDe.Range("H4").Value = "2016-08-24 22:46:12"
$wsDeDe.Range("K4").Value = "2016-08-24 22:46:32"
